$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.08632499999999999
$ws.Range("H2").Value = 0.258975
$ws.Range("I2").Value = 0.04465318711422561
$ws.Range("J2").Value = 0.0446531871142256
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.325274
$ws.Range("N2").Value = 0.975822
$ws.Range("O2").Value = 0.4689956999283401
$ws.Range("P2").Value = 0.4689956999283402
$ws.Range("Q2").Value = 0.02807927805
$ws.Range("R2").Value = 0.25271350245
$ws.Range("S2").Value = 0.02094215274466737
$ws.Range("T2").Value = 0.02094215274466737
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.08632499999999999
$ws.Range("H3").Value = 0.258975
$ws.Range("I3").Value = 0.04465318711422561
$ws.Range("J3").Value = 0.0446531871142256
$ws.Range("M3").Value = 0.3682803333333333
$ws.Range("N3").Value = 1.104841
$ws.Range("O3").Value = 0.5310043000716598
$ws.Range("P3").Value = 0.5310043000716599
$ws.Range("Q3").Value = 0.03179179977499999
$ws.Range("R3").Value = 0.2861261979749999
$ws.Range("S3").Value = 0.02371103436955822
$ws.Range("T3").Value = 0.02371103436955823
$ws.Range("I4").Value = 0.2144309927861761
$ws.Range("J4").Value = 0.2144309927861761
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.325274
$ws.Range("N4").Value = 0.975822
$ws.Range("O4").Value = 0.4689956999283401
$ws.Range("P4").Value = 0.4689956999283402
$ws.Range("Q4").Value = 0.13484071033
$ws.Range("R4").Value = 1.21356639297
$ws.Range("S4").Value = 0.1005672135480815
$ws.Range("T4").Value = 0.1005672135480815
$ws.Range("I5").Value = 0.2144309927861761
$ws.Range("J5").Value = 0.2144309927861761
$ws.Range("M5").Value = 0.3682803333333333
$ws.Range("N5").Value = 1.104841
$ws.Range("O5").Value = 0.5310043000716598
$ws.Range("P5").Value = 0.5310043000716599
$ws.Range("Q5").Value = 0.1526687707816667
$ws.Range("R5").Value = 1.374018937035
$ws.Range("S5").Value = 0.1138637792380946
$ws.Range("T5").Value = 0.1138637792380946
$ws.Range("G6").Value = 1.275124666666666
$ws.Range("H6").Value = 3.825374
$ws.Range("I6").Value = 0.6595815851101212
$ws.Range("J6").Value = 0.6595815851101212
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.325274
$ws.Range("N6").Value = 0.975822
$ws.Range("O6").Value = 0.4689956999283401
$ws.Range("P6").Value = 0.4689956999283402
$ws.Range("Q6").Value = 0.4147649008253333
$ws.Range("R6").Value = 3.732884107428
$ws.Range("S6").Value = 0.3093409271685653
$ws.Range("T6").Value = 0.3093409271685654
$ws.Range("G7").Value = 1.275124666666666
$ws.Range("H7").Value = 3.825374
$ws.Range("I7").Value = 0.6595815851101212
$ws.Range("J7").Value = 0.6595815851101212
$ws.Range("M7").Value = 0.3682803333333333
$ws.Range("N7").Value = 1.104841
$ws.Range("O7").Value = 0.5310043000716598
$ws.Range("P7").Value = 0.5310043000716599
$ws.Range("Q7").Value = 0.4696033372815555
$ws.Range("R7").Value = 4.226430035533999
$ws.Range("S7").Value = 0.3502406579415558
$ws.Range("T7").Value = 0.3502406579415559
$ws.Range("G8").Value = 0.157238
$ws.Range("H8").Value = 0.471714
$ws.Range("I8").Value = 0.08133423498947705
$ws.Range("J8").Value = 0.08133423498947705
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.325274
$ws.Range("N8").Value = 0.975822
$ws.Range("O8").Value = 0.4689956999283401
$ws.Range("P8").Value = 0.4689956999283402
$ws.Range("Q8").Value = 0.051145433212
$ws.Range("R8").Value = 0.460308898908
$ws.Range("S8").Value = 0.03814540646702588
$ws.Range("T8").Value = 0.03814540646702588
$ws.Range("G9").Value = 0.157238
$ws.Range("H9").Value = 0.471714
$ws.Range("I9").Value = 0.08133423498947705
$ws.Range("J9").Value = 0.08133423498947705
$ws.Range("M9").Value = 0.3682803333333333
$ws.Range("N9").Value = 1.104841
$ws.Range("O9").Value = 0.5310043000716598
$ws.Range("P9").Value = 0.5310043000716599
$ws.Range("Q9").Value = 0.05790766305266666
$ws.Range("R9").Value = 0.521168967474
$ws.Range("S9").Value = 0.04318882852245116
$ws.Range("T9").Value = 0.04318882852245117
